$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("A2").Value = 2993253
$ws.Range("D2").Value = 223
$ws.Range("D3").Value = 223

# Reset the view: scroll back to top-left and select a single cell (D4)
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("D4").Select()
